$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.186.43'
$ws.Range("E2").Value = '  +4.57%  '
$ws.Range("D3").Value = '3.243.64'
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.26'
$ws.Range("E5").Value = '  +2.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.49'
$ws.Range("E6").Value = '  +5.58%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  -1.25%  '
$ws.Range("D9").Value = '3.241.81'
$ws.Range("E9").Value = '  +2.04%  '
$ws.Range("E10").Value = '  +4.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.74'
$ws.Range("E11").Value = '  +1.82%  '
$ws.Range("E12").Value = '  +4.47%  '
$ws.Range("D13").Value = '3.804.78'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.91'
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").Value = '67.115.65'
$ws.Range("E16").Value = '  +4.47%  '
$ws.Range("E17").Value = '  +2.80%  '
$ws.Range("D18").Value = '3.245.31'
$ws.Range("E18").Value = '  +2.18%  '
$ws.Range("E19").Value = '  +2.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.31'
$ws.Range("E20").Value = '  +2.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.86'
$ws.Range("E21").Value = '  +5.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.60'
$ws.Range("E22").Value = '  +5.91%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.07'
$ws.Range("E24").Value = '  +3.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.509'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '3.381.10'
$ws.Range("E26").Value = '  +2.12%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.87'
$ws.Range("E28").Value = '  +2.82%  '
$ws.Range("E29").Value = '  +1.89%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("E31").Value = '  +3.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.63'
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.54'
$ws.Range("E33").Value = '  +2.54%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.28'
$ws.Range("E35").Value = '  +6.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.81'
$ws.Range("E36").Value = '  +2.58%  '
$ws.Range("E37").Value = '  +4.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.62'
$ws.Range("E38").Value = '  +4.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.859'
$ws.Range("E39").Value = '  +5.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.85'
$ws.Range("E40").Value = '  +9.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.85'
$ws.Range("E41").Value = '  +14.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.80'
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("E43").Value = '  +4.76%  '
$ws.Range("D44").Value = '2.765.23'
$ws.Range("E44").Value = '  +5.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '357.64'
$ws.Range("E45").Value = '  +12.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.40'
$ws.Range("E46").Value = '  +5.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.79'
$ws.Range("E47").Value = '  +8.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.41'
$ws.Range("E48").Value = '  +2.61%  '
$ws.Range("E49").Value = '  +3.89%  '
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.103'
$ws.Range("E51").Value = '  +1.49%  '
